# Update Excel data - 2024-11-22 06:06:33
# Refreshes the live crypto-market snapshot across all three sheets:
#   - "Top 50 Cryptocurrencies": per-coin price / market-cap / volume / change
#     (rows 22-23 also swap: Polkadot now ranks above Pepe)
#   - "Top 5 by Market Cap": refreshed market caps for the top 5 coins
#   - "Summary": refreshed derived metrics (average price, highest/lowest 24h change)

$wb = $excel.ActiveWorkbook

# Writes a value as literal text even when it looks like a number/currency
# (e.g. "$4354.78"), so Excel's COM auto-conversion doesn't turn it into a
# numeric cell. Restores "General" formatting afterwards so the cell's style
# index is left untouched (matches cells that were never specially formatted).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet 1: "Top 50 Cryptocurrencies" -- columns A:Name B:Symbol C:Price
#          D:MarketCap E:Volume24h F:Change24h, data rows 2..51
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$sheet1Data = @(
    ('Bitcoin','btc',98949,1957312348128,105962937556,1.8296),
    ('Ethereum','eth',3385.16,407528946484,57176037571,9.096640000000001),
    ('Tether','usdt',1.001,130860870762,126294469787,0.12279),
    ('Solana','sol',261.06,123886385324,14815108915,9.37628),
    ('BNB','bnb',633.28,92366163407,2380817573,4.06307),
    ('XRP','xrp',1.4,79935603727,18250929448,26.4852),
    ('Dogecoin','doge',0.396128,58217902823,9750838152,3.58956),
    ('USDC','usdc',1,38321674118,12198920959,0.01437),
    ('Lido Staked Ether','steth',3386.24,33159664554,143198614,8.99356),
    ('Cardano','ada',0.882772,31541459969,3550138063,12.56423),
    ('TRON','trx',0.20009,17277520321,1062158397,1.76703),
    ('Avalanche','avax',36.41,14883697914,1047195719,7.4404),
    ('Shiba Inu','shib',0.000025,14724458250,1598057036,3.9696),
    ('Wrapped stETH','wsteth',4005.46,14482650687,161361117,8.896800000000001),
    ('Wrapped Bitcoin','wbtc',98714,14414801790,837289225,1.93306),
    ('Toncoin','ton',5.55,14140621430,624413906,3.71466),
    ('Sui','sui',3.61,10262277354,2046611540,0.5314),
    ('Bitcoin Cash','bch',498.98,9882578564,1817006727,-3.37977),
    ('WETH','weth',3383.73,9695241516,2193400238,8.92578),
    ('Chainlink','link',15.32,9597641840,1256233231,5.70752),
    ('Polkadot','dot',6.24,8977044161,839252580,10.40347),
    ('Pepe','pepe',0.00002133,8974266732,6734296727,10.19788),
    ('Stellar','xlm',0.284633,8528344379,2313085743,20.34539),
    ('LEO Token','leo',8.800000000000001,8138880715,3421347,3.3504),
    ('NEAR Protocol','near',5.82,7080828321,1010209673,5.71874),
    ('Litecoin','ltc',91.03,6847455412,1385604498,5.40234),
    ('Aptos','apt',12.14,6472209222,856252554,4.43628),
    ('Wrapped eETH','weeth',3566.33,6211334766,105658696,9.08948),
    ('Uniswap','uni',9.390000000000001,5640414358,870529475,6.54853),
    ('Cronos','cro',0.202831,5533748262,136788438,16.49178),
    ('USDS','usds',0.999576,5228866032,16084967,0.13728),
    ('Hedera','hbar',0.136113,5156068985,921621803,10.399),
    ('Internet Computer','icp',9.68,4584804726,275557543,7.56923),
    ('Ethereum Classic','etc',28.13,4210412676,870064263,6.04145),
    ('Bonk','bonk',0.00005225,3919437831,1639904482,1.73519),
    ('Kaspa','kas',0.152416,3840028133,151628290,1.53692),
    ('Render','render',7.41,3835215289,431243780,1.12561),
    ('POL (ex-MATIC)','pol',0.469936,3758408261,496353927,8.170859999999999),
    ('Bittensor','tao',504.32,3722424094,281336572,3.64941),
    ('Ethena USDe','usde',1.002,3689859605,224518933,0.13789),
    ('WhiteBIT Coin','wbt',24.88,3575145979,44086915,2.94399),
    ('Dai','dai',1,3441676073,156466688,-0.00268),
    ('dogwifhat','wif',3.39,3382728022,1286539243,5.68177),
    ('MANTRA','om',3.78,3381419859,302623861,5.24868),
    ('Artificial Superintelligence Alliance','fet',1.28,3344707679,484485647,3.49845),
    ('Arbitrum','arb',0.793462,3244693440,1673139119,15.0305),
    ('Monero','xmr',161.26,2974960015,86422219,-0.45054),
    ('Stacks','stx',1.97,2957274297,355653430,2.91976),
    ('Mantle','mnt',0.84516,2841679667,183193974,16.02787),
    ('Filecoin','fil',4.71,2830947251,576676347,8.232100000000001)
)

for ($i = 0; $i -lt $sheet1Data.Count; $i++) {
    $row = $i + 2
    $vals = $sheet1Data[$i]
    $ws1.Cells.Item($row, 1).Value = $vals[0]
    $ws1.Cells.Item($row, 2).Value = $vals[1]
    $ws1.Cells.Item($row, 3).Value = $vals[2]
    $ws1.Cells.Item($row, 4).Value = $vals[3]
    $ws1.Cells.Item($row, 5).Value = $vals[4]
    $ws1.Cells.Item($row, 6).Value = $vals[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Top 5 by Market Cap" -- columns A:Name B:MarketCap, rows 2..6
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")

$sheet2Data = @(
    ('Bitcoin',1957312348128),
    ('Ethereum',407528946484),
    ('Tether',130860870762),
    ('Solana',123886385324),
    ('BNB',92366163407)
)

for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $row = $i + 2
    $vals = $sheet2Data[$i]
    $ws2.Cells.Item($row, 1).Value = $vals[0]
    $ws2.Cells.Item($row, 2).Value = $vals[1]
}

# ---------------------------------------------------------------------------
# Sheet 3: "Summary" -- columns A:Metric B:Value, rows 2..4
# (B2 holds a dollar-formatted string, e.g. "$4354.78", which must stay text)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Summary")

$sheet3Data = @(
    ('Average Price','$4354.78'),
    ('Highest 24h Change','XRP (26.49%)'),
    ('Lowest 24h Change','Bitcoin Cash (-3.38%)')
)

for ($i = 0; $i -lt $sheet3Data.Count; $i++) {
    $row = $i + 2
    $vals = $sheet3Data[$i]
    $ws3.Cells.Item($row, 1).Value = $vals[0]
    Set-TextValue $ws3.Cells.Item($row, 2) $vals[1]
}
